$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 ("Save"), copying the same formatting (style) as G1 ("sum")
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data cells H2, H3 with value 0
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
